# Netflix.xlsx update — "update year movie and net"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Existing "연간" (annual) table: A17/A18 switch from the 2-dp
#    paren format to an integer paren format (still whole years).
# ---------------------------------------------------------------
$ws.Range("A17:A18").NumberFormat = "0_);[Red]\(0\)"

# ---------------------------------------------------------------
# 2. New row 19 — Q1 2021 style annual year addition
# ---------------------------------------------------------------
$ws.Range("A19").Value = 2021

# ---------------------------------------------------------------
# 3. New data dropped next to the existing "넷플릭스 컨텐츠" block
# ---------------------------------------------------------------
$ws.Range("E24").Value = 1948
$ws.Range("E25").Value = 1899
$ws.Range("F27").Value = 1913989080068

# ---------------------------------------------------------------
# 4. New header row 30 (mirrors row 1 / row 16 style, reuses the
#    same 2-dp paren number format / vertical-center style as the
#    rest of the header rows)
# ---------------------------------------------------------------
$ws.Range("A30").Value = "Year "
$ws.Range("B30").Value = "Asia_sub"
$ws.Range("C30").Value = "Sub_growth "
$ws.Range("D30").Value = "Asia_revenue "
$ws.Range("E30").Value = "Net_revenue_growth"
$ws.Range("F30").Value = "Movies"
$ws.Range("G30").Value = "Movie_revenue "
$ws.Range("H30").Value = "Movie_growth"
$ws.Range("I30").Value = "Audience"
$ws.Range("J30").Value = "Audience_growth"
$ws.Range("K30").Value = "Net_contents"
$ws.Range("A30:K30").NumberFormat = "0.00_);[Red]\(0.00\)"

# ---------------------------------------------------------------
# 5. Row 31 / 32 — the new "Year / Asia_sub / ..." summary table
# ---------------------------------------------------------------
$ws.Range("A31").Value = 2019
$ws.Range("A32").Value = 2020
$ws.Range("A31:A32").NumberFormat = "0_);[Red]\(0\)"

$ws.Range("B31").Value = 55.8
$ws.Range("B32").Value = 91.32
$ws.Range("C31").Value = ""
$ws.Range("C32").Formula = "=(B32-B31)/B31*100"
$ws.Range("D31").Value = 1469
$ws.Range("D32").Value = 2373
$ws.Range("E31").Value = ""
$ws.Range("E32").Formula = "=(D32-D31)/D31*100"
$ws.Range("B31:E32").NumberFormat = "0.00_);[Red]\(0.00\)"

$ws.Range("B24").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F32").PasteSpecial(-4122)
$ws.Range("F31").Value = 5769
$ws.Range("F32").Value = 5838

$ws.Range("G31").Value = 1913989.08
$ws.Range("G32").Value = 510375.22219100001
$ws.Range("H32").Value = -73
$ws.Range("J32").Value = -74
$ws.Range("G31:G32").NumberFormat = "0.00_ "
$ws.Range("H31").NumberFormat = "0_ "
$ws.Range("H32").NumberFormat = "0.00_ "
$ws.Range("J32").NumberFormat = "0.00_ "

$ws.Range("I31").Value = 226.678777
$ws.Range("I32").Value = 59.524093000000001
$ws.Range("I31:I32").NumberFormat = "0.00;[Red]0.00"

$ws.Range("K31").Value = 1948
$ws.Range("K32").Value = 1899

# ---------------------------------------------------------------
# 6. Row 35 / 36 — same figures, flattened into a contiguous block
# ---------------------------------------------------------------
$ws.Range("A35").Value = 2019
$ws.Range("A36").Value = 2020
$ws.Range("A35:A36").NumberFormat = "0_);[Red]\(0\)"

$ws.Range("B35").Value = 55.8
$ws.Range("B36").Value = 91.32
$ws.Range("C35").Value = 1469
$ws.Range("C36").Value = 2373
$ws.Range("B35:C36").NumberFormat = "0.00_);[Red]\(0.00\)"

$ws.Range("B24").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D35").Value = 5769
$ws.Range("D36").Value = 5838

$ws.Range("E35").Value = 1913989.08
$ws.Range("E36").Value = 510375.22219100001
$ws.Range("E35:E36").NumberFormat = "0.00_ "

$ws.Range("F35").Value = 226.678777
$ws.Range("F36").Value = 59.524093000000001
$ws.Range("F35:F36").NumberFormat = "0.00;[Red]0.00"

$ws.Range("G35").Value = 1948
$ws.Range("G36").Value = 1899

# ---------------------------------------------------------------
# 7. Restore selection roughly where the author left it
# ---------------------------------------------------------------
$ws.Range("I22").Select()
